# The Keywords column (C) held a few values that needed cleanup:
#  - numeric-looking keyword text that should match other approximations
#  - slash-separated keyword lists changed to space-separated lists
#  - a couple of keyword lists trimmed down
#
# A leading apostrophe forces Excel to store these as literal text (not
# numbers), matching how the rest of the Keywords column is stored; the
# ClearFormats() afterwards drops the resulting quote-prefix cell style so
# only the cell's value/type changes, not its formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    if ($text -match '^-?[0-9]') {
        $range.Value = "'" + $text
        $range.ClearFormats()
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range("C12") "3.14"
Set-TextValue $ws.Range("C14") "HTML"
Set-TextValue $ws.Range("C18") "attention careful"
Set-TextValue $ws.Range("C20") "Rendering Navigation Interaction"
Set-TextValue $ws.Range("C21") "Cacti Cactus"
Set-TextValue $ws.Range("C22") "nucleus"
Set-TextValue $ws.Range("C23") "299.0"
